$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new rows first, working from the bottom of the sheet upward
#     so each Rows().Insert() call can use the still-valid original row
#     numbers for everything above it. ---

# "Block Power Converter" section: room for 1 new part (E36SR12004), right
# below the section header, before the existing blank separator row.
$ws.Rows("21").Insert()

# "Blocking Diodes" section: room for 2 new parts (STPS5L60, MBR360GOS),
# right below the section header, before the existing blank separator row.
$ws.Rows("19:20").Insert()

# "TVS Diode" section: room for 2 new parts (SMAJ16CA, SMBJ60CA) right below
# the section header, consuming one of the two existing blank rows.
$ws.Rows("16").Insert()

# --- Now fill in the new values. ---

# Newly-calibrated FC current sensor part, under "Block Power Converter".
$ws.Range("A24").Value = "E36SR12004"
$ws.Range("C24").Value = 1

# New safety parts under "TVS Diode" and "Blocking Diodes", entered as two
# passes (first new part of each section, then the second of each).
$ws.Range("A16").Value = "SMAJ16CA"
$ws.Range("C16").Value = 1
$ws.Range("A20").Value = "STPS5L60"
$ws.Range("C20").Value = 3
$ws.Range("A17").Value = "SMBJ60CA"
$ws.Range("C17").Value = 1
$ws.Range("A21").Value = "MBR360GOS"
$ws.Range("C21").Value = 1

# Capacitors section: "3.3 uF" now has a quantity of 2.
$ws.Range("C5").Value = 2

# "12 to 5 regulator" section: LM2576D2TR4-5GOSCT-ND now has a quantity of 1.
$ws.Range("C27").Value = 1

# Park the selection on the row past the end of the list, matching where the
# editor left off after making these changes.
$ws.Range("A32").Select()
